# Applies the data-table changes described by the commit:
#   "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The "Periodo Mora" list (rows 16-26, column E) gets reversed, and the
# corresponding "Valor Mora" (column F) value that used to sit on the
# first period (1902) now sits on the last period (1912) and vice versa.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 16
$lastRow  = 26

# Capture current (pre-edit) E/F values for the block so we can reverse them.
$periods = @()
$values  = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $periods += , $ws.Range("E$r").Value2
    $values  += , $ws.Range("F$r").Value2
}

# Reverse the period labels; keep the Valor Mora values aligned to their
# original row position (only the two rows whose value actually differs -
# first/last - end up changing).
$revPeriods = @($periods[($periods.Length - 1)..0])

for ($i = 0; $i -lt $periods.Length; $i++) {
    $r = $firstRow + $i
    $ws.Range("E$r").Value = $revPeriods[$i]
}

# Swap the Valor Mora figures for the first and last rows of the block
# (30400 <-> 38000), matching the new period order.
$ws.Range("F$firstRow").Value = $values[$values.Length - 1]
$ws.Range("F$lastRow").Value  = $values[0]
